$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# --- Step 1: rewrite the "-lastlogin: String" paragraph, dropping the
#     "_GoBack" bookmark that currently lives there (it moves to the new
#     row below, tracking the most recent edit point like real Word does).
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*-lastlogin: String*") {
        $lastLoginXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6ACD12F9" w14:textId="5C66B86D" w:rsidR="007138E0" w:rsidRDefault="002111F9" w:rsidP="00BD2C13"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>-lastlogin: String</w:t></w:r></w:p>
'@
        $para.Range.InsertXML($lastLoginXml)
        break
    }
}

# --- Step 2: append a new table row holding the staff-service method
#     signatures, with the "_GoBack" bookmark re-anchored in the
#     "approveBeneficiary();" paragraph (the point of the last edit).
$newRow = $table.Rows.Add()
$newCell = $newRow.Cells.Item(1)

$newRowXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>getID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>getName</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>getDOB</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>getDepartment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>getAddress</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>getDOJ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>getLastLogin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>getRequestedBeneficiary</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>approveBeneficiary</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>getRequestedChequeBook</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>approveChequeBook</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>changePassword</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>+</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>logout(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>);</w:t></w:r></w:p>
'@
$newCell.Range.InsertXML($newRowXml)
